# Updates the cryptocurrency price table (rows 2-51) on Sheet1 to the latest
# scrape: a new coin (WrappedliquidstakedEther2.0) enters the top of the
# ranking around row 17, pushing lower-ranked coins down by one row (and
# Cronos drops off the bottom), plus refreshed Price/Volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="29.425.42"; E="  -0.59%  " },
    @{ Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="1.848.19"; E="  -0.61%  " },
    @{ Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="0.9985"; E="  +0.00%  " },
    @{ Row=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="240.71"; E="  -0.80%  " },
    @{ Row=6; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.6317"; E="  -0.52%  " },
    @{ Row=7; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="0.9998"; E="  +0.06%  " },
    @{ Row=8; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.07561"; E="  -0.57%  " },
    @{ Row=9; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.2953"; E="  -1.60%  " },
    @{ Row=10; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="24.56"; E="  -0.56%  " },
    @{ Row=11; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.07696"; E="  -0.64%  " },
    @{ Row=12; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.852.64"; E="  -0.30%  " },
    @{ Row=13; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="4.991"; E="  -1.00%  " },
    @{ Row=14; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.6858"; E="  -1.60%  " },
    @{ Row=15; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.00001006"; E="  +0.69%  " },
    @{ Row=16; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="83.11"; E="  -0.90%  " },
    @{ Row=17; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.107.46"; E="  -0.25%  " },
    @{ Row=18; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="6.131"; E="  -2.44%  " },
    @{ Row=19; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="29.439.26"; E="  -0.58%  " },
    @{ Row=20; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="228.58"; E="  -2.93%  " },
    @{ Row=21; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="12.49"; E="  -0.99%  " },
    @{ Row=22; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="0.9995"; E="  -0.07%  " },
    @{ Row=23; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="7.540"; E="  -1.96%  " },
    @{ Row=24; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.000"; E="  +0.08%  " },
    @{ Row=25; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="157.08"; E="  +0.51%  " },
    @{ Row=26; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.1394"; E="  -0.71%  " },
    @{ Row=27; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="8.372"; E="  -1.55%  " },
    @{ Row=28; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="17.68"; E="  -0.66%  " },
    @{ Row=29; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="1.468"; E="  -0.64%  " },
    @{ Row=30; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="1.266"; E="  +0.25%  " },
    @{ Row=31; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.05731"; E="  -1.74%  " },
    @{ Row=32; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.122"; E="  -0.53%  " },
    @{ Row=33; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="4.026"; E="  -0.43%  " },
    @{ Row=34; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.846"; E="  -3.54%  " },
    @{ Row=35; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.155"; E="  -1.53%  " },
    @{ Row=36; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.7145"; E="  -1.22%  " },
    @{ Row=37; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.587"; E="  +0.05%  " },
    @{ Row=38; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="1.249.03"; E="  -0.65%  " },
    @{ Row=39; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01808"; E="  -0.06%  " },
    @{ Row=40; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.776"; E="  -1.07%  " },
    @{ Row=41; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="0.9108"; E="  +0.13%  " },
    @{ Row=42; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="6.171"; E="  +0.40%  " },
    @{ Row=43; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="1.000"; E="  +0.09%  " },
    @{ Row=44; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="101.85"; E="  +0.24%  " },
    @{ Row=45; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="66.09"; E="  -3.63%  " },
    @{ Row=46; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="7.090"; E="  -3.74%  " },
    @{ Row=47; B="BabyDogeCoin"; C="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D="0.00000000118"; E="  -0.48%  " },
    @{ Row=48; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.4024"; E="  -1.07%  " },
    @{ Row=49; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="9.080"; E="  -1.02%  " },
    @{ Row=50; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="1.687"; E="  -2.01%  " },
    @{ Row=51; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1120"; E="  -0.26%  " }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    # Price strings such as "4.991" or "0.1120" parse as genuine numbers in
    # Excel, which would silently drop the meaningful trailing/leading zeros
    # (the source site formats these as text). Prefix with an apostrophe so
    # they stay literal text, same as typing them in by hand.
    if ($r.D -match '^[0-9]+(\.[0-9]+)?$') {
        $ws.Cells.Item($r.Row, 4).Value = "'" + $r.D
    } else {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
